# Apply the "data updated on Aug.09" edit:
#  - fix the stray thick left-border on F29:F33 (make them match the
#    normal/no-border style used elsewhere in the table)
#  - de-duplicate the D/E/H/I borderless style used on rows 31-33
#    (same visual result, just collapses onto the same style as the
#    rest of the sheet)
#  - rows 29 & 30 take the slightly shorter 13.8pt row height already
#    used by rows 31-33
#  - append a new data row (row 34) for 2022-08-08 (serial 44781) with
#    all-zero counts, formatted like the rows above it

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fix anomalous thick-left-border cells in F29:F33, and de-dupe the
#     borderless style used by D31:E33 / H31:I33 so everything lines up
#     with the rest of the table (no visible border in either case) ---
$ws.Range("F29:F33").Borders.LineStyle = -4142
$ws.Range("D31:E33").Borders.LineStyle = -4142
$ws.Range("H31:I33").Borders.LineStyle = -4142

# --- rows 29 & 30 get the shorter row height used by 31-33 ---
$ws.Rows.Item(29).RowHeight = 13.8
$ws.Rows.Item(30).RowHeight = 13.8

# --- add new row 34 (2022-08-08, serial 44781) ---
# copy formatting from row 33 (already fixed above) so fonts/number
# formats/row height match the rest of the table, then fill in values
$ws.Range("A33:I33").Copy()
$ws.Range("A34:I34").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Rows.Item(34).RowHeight = 13.8

$ws.Cells.Item(34, 1).Value2 = 44781
$ws.Cells.Item(34, 2).Value2 = 0
$ws.Cells.Item(34, 3).Value2 = 0
$ws.Cells.Item(34, 4).Value2 = 0
$ws.Cells.Item(34, 5).Value2 = 0
$ws.Cells.Item(34, 6).Value2 = 0
$ws.Cells.Item(34, 7).Value2 = 0
$ws.Cells.Item(34, 8).Value2 = 0
$ws.Cells.Item(34, 9).Value2 = 0

# F34 should match the borderless style (no thick left border)
$ws.Range("F34").Borders.LineStyle = -4142

# --- update the view so the new last rows are visible/selected ---
$ws.Range("G39").Select() | Out-Null

Write-Host "Edit applied"
